# Insert a new row for "llama3_70b_instruct" after the header rows,
# right before "llama3_8b_instruct" (currently row 5), shifting the
# rest of the table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 5; existing row 5 and below shift to row 6 and below.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the llama3_70b_instruct metrics.
$ws.Range("A5").Value = "llama3_70b_instruct"
$ws.Range("B5").Value = 2296
$ws.Range("C5").Value = 1941
$ws.Range("D5").Value = 1959
$ws.Range("E5").Value = 1822
$ws.Range("F5").Value = 272
$ws.Range("G5").Value = 222
$ws.Range("H5").Value = 185
$ws.Range("I5").Value = 1
